$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (22 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1961252.2
$ws.Range("J17").Value = 2000464
$ws.Range("L17").Value = 6001392
$ws.Range("N17").Value = -6001728
$ws.Range("H32").Value = 523.1667
$ws.Range("I32").Value = 350
$ws.Range("J32").Value = 696.3333
$ws.Range("K32").Value = 350
$ws.Range("L32").Value = 696.3333
$ws.Range("M32").Value = -24
$ws.Range("N32").Value = -1348.3333
$ws.Range("H129").Value = 162126.05
$ws.Range("J129").Value = 170325.84
$ws.Range("L129").Value = 510977.52
$ws.Range("N129").Value = -520977.52
$ws.Range("H138").Value = 1401.8387
$ws.Range("I138").Value = 548.2381
$ws.Range("J138").Value = 3194.4
$ws.Range("K138").Value = 1644.7143
$ws.Range("L138").Value = 9583.200000000001
$ws.Range("M138").Value = 3495.2857
$ws.Range("N138").Value = -19863.2

# --- Sheet: ARM (50 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2469.3818
$ws.Range("I32").Value = 1820.0834
$ws.Range("J32").Value = 6921.7144
$ws.Range("K32").Value = 1820.0834
$ws.Range("L32").Value = 6921.7144
$ws.Range("M32").Value = -1533.0834
$ws.Range("N32").Value = -7495.7144
$ws.Range("H45").Value = 3298.3845
$ws.Range("I45").Value = 3523.1428
$ws.Range("K45").Value = 3523.1428
$ws.Range("M45").Value = -3146.1428
$ws.Range("H61").Value = 3222.8635
$ws.Range("I61").Value = 2218.9375
$ws.Range("J61").Value = 5900
$ws.Range("K61").Value = 2218.9375
$ws.Range("L61").Value = 5900
$ws.Range("M61").Value = -2006.9375
$ws.Range("N61").Value = -6324
$ws.Range("H74").Value = 2615.3333
$ws.Range("I74").Value = 2777.5
$ws.Range("J74").Value = 1966.6666
$ws.Range("K74").Value = 2777.5
$ws.Range("L74").Value = 1966.6666
$ws.Range("M74").Value = -1903.5
$ws.Range("N74").Value = -3714.6666
$ws.Range("H77").Value = 2615.3333
$ws.Range("I77").Value = 2777.5
$ws.Range("J77").Value = 1966.6666
$ws.Range("K77").Value = 13887.5
$ws.Range("L77").Value = 9833.333000000001
$ws.Range("M77").Value = -9519.5
$ws.Range("N77").Value = -18569.333
$ws.Range("H102").Value = 3098.6
$ws.Range("I102").Value = 1162.5
$ws.Range("K102").Value = 1162.5
$ws.Range("M102").Value = 459.5
$ws.Range("H122").Value = 1459.6818
$ws.Range("I122").Value = 1487.6
$ws.Range("J122").Value = 1399.8572
$ws.Range("K122").Value = 4462.799999999999
$ws.Range("L122").Value = 4199.571599999999
$ws.Range("M122").Value = -2012.799999999999
$ws.Range("N122").Value = -9099.571599999999
$ws.Range("H136").Value = 3222.8635
$ws.Range("I136").Value = 2218.9375
$ws.Range("J136").Value = 5900
$ws.Range("K136").Value = 6656.8125
$ws.Range("L136").Value = 17700
$ws.Range("M136").Value = -4106.8125
$ws.Range("N136").Value = -22800

# --- Sheet: BSM (22 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2304.762
$ws.Range("I86").Value = 1887.4117
$ws.Range("J86").Value = 4078.5
$ws.Range("K86").Value = 1887.4117
$ws.Range("L86").Value = 4078.5
$ws.Range("M86").Value = -764.4117000000001
$ws.Range("N86").Value = -6324.5
$ws.Range("H89").Value = 2304.762
$ws.Range("I89").Value = 1887.4117
$ws.Range("J89").Value = 4078.5
$ws.Range("K89").Value = 9437.058500000001
$ws.Range("L89").Value = 20392.5
$ws.Range("M89").Value = -3821.058500000001
$ws.Range("N89").Value = -31624.5
$ws.Range("H99").Value = 2225.5715
$ws.Range("I99").Value = 1766.5
$ws.Range("K99").Value = 1766.5
$ws.Range("M99").Value = -268.5
$ws.Range("H134").Value = 3082.7778
$ws.Range("I134").Value = 3217.4
$ws.Range("K134").Value = 9652.200000000001
$ws.Range("M134").Value = -7117.200000000001

# --- Sheet: CRP (34 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 849
$ws.Range("I16").Value = 807.625
$ws.Range("J16").Value = 959.3333
$ws.Range("K16").Value = 807.625
$ws.Range("L16").Value = 959.3333
$ws.Range("M16").Value = -520.625
$ws.Range("N16").Value = -1533.3333
$ws.Range("H31").Value = 9321.976000000001
$ws.Range("I31").Value = 10451.529
$ws.Range("K31").Value = 10451.529
$ws.Range("M31").Value = -10156.529
$ws.Range("H34").Value = 9321.976000000001
$ws.Range("I34").Value = 10451.529
$ws.Range("K34").Value = 10451.529
$ws.Range("M34").Value = -10249.529
$ws.Range("H58").Value = 19654.814
$ws.Range("I58").Value = 1173.8182
$ws.Range("K58").Value = 1173.8182
$ws.Range("M58").Value = -970.8181999999999
$ws.Range("H113").Value = 849
$ws.Range("I113").Value = 807.625
$ws.Range("J113").Value = 959.3333
$ws.Range("K113").Value = 807.625
$ws.Range("L113").Value = 959.3333
$ws.Range("M113").Value = 1362.375
$ws.Range("N113").Value = -5299.3333
$ws.Range("H125").Value = 13000
$ws.Range("J125").Value = 13000
$ws.Range("L125").Value = 13000
$ws.Range("N125").Value = -17920
$ws.Range("H136").Value = 19654.814
$ws.Range("I136").Value = 1173.8182
$ws.Range("K136").Value = 3521.4546
$ws.Range("M136").Value = -971.4546

# --- Sheet: CUL (15 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 898.5714
$ws.Range("I92").Value = 533.3333
$ws.Range("J92").Value = 1172.5
$ws.Range("K92").Value = 1599.9999
$ws.Range("L92").Value = 3517.5
$ws.Range("M92").Value = -351.9999
$ws.Range("N92").Value = -6013.5
$ws.Range("H131").Value = 755.84
$ws.Range("J131").Value = 769.6288500000001
$ws.Range("L131").Value = 2308.88655
$ws.Range("N131").Value = -12388.88655
$ws.Range("H133").Value = 3920
$ws.Range("I133").Value = 3840
$ws.Range("K133").Value = 11520
$ws.Range("M133").Value = -6460

# --- Sheet: GSM (7 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1740.5
$ws.Range("I122").Value = 1719.8
$ws.Range("J122").Value = 1761.2
$ws.Range("K122").Value = 5159.4
$ws.Range("L122").Value = 5283.6
$ws.Range("M122").Value = -2709.4
$ws.Range("N122").Value = -10183.6

# --- Sheet: LTW (25 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6187.375
$ws.Range("I40").Value = 5499.8
$ws.Range("J40").Value = 7333.3335
$ws.Range("K40").Value = 5499.8
$ws.Range("L40").Value = 7333.3335
$ws.Range("M40").Value = -5363.8
$ws.Range("N40").Value = -7605.3335
$ws.Range("H46").Value = 2382.8572
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 2670
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2670
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -3046
$ws.Range("H132").Value = 3371.8572
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 4120.6
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 12361.8
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -17421.8
$ws.Range("H136").Value = 500000.5
$ws.Range("I136").Value = 500000.5
$ws.Range("K136").Value = 1500001.5
$ws.Range("M136").Value = -1497451.5

# --- Sheet: WVR (40 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2066.6667
$ws.Range("I81").Value = 2066.6667
$ws.Range("K81").Value = 4133.3334
$ws.Range("M81").Value = -3072.3334
$ws.Range("H84").Value = 2066.6667
$ws.Range("I84").Value = 2066.6667
$ws.Range("K84").Value = 20666.667
$ws.Range("M84").Value = -15362.667
$ws.Range("H107").Value = 3497788.2
$ws.Range("I107").Value = 440
$ws.Range("J107").Value = 5683631
$ws.Range("K107").Value = 1320
$ws.Range("L107").Value = 17050893
$ws.Range("M107").Value = 600
$ws.Range("N107").Value = -17054733
$ws.Range("H122").Value = 1598.9565
$ws.Range("I122").Value = 1464.8948
$ws.Range("K122").Value = 4394.6844
$ws.Range("M122").Value = -1944.6844
$ws.Range("H126").Value = 1932.25
$ws.Range("I126").Value = 1051.6
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 3154.8
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -684.7999999999997
$ws.Range("N126").Value = -15140
$ws.Range("H132").Value = 3296.158
$ws.Range("I132").Value = 2971.6155
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 8914.8465
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -6384.8465
$ws.Range("N132").Value = -17057.9999
$ws.Range("H136").Value = 33334650
$ws.Range("I136").Value = 50000970
$ws.Range("J136").Value = 2010.4
$ws.Range("K136").Value = 150002910
$ws.Range("L136").Value = 6031.200000000001
$ws.Range("M136").Value = -150000360
$ws.Range("N136").Value = -11131.2

Write-Output "Applied $([int]215) cell updates"